# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Lane Late - Primera / Segunda) for
# Femacal de La Calera - Naranja right after the existing row 441 block,
# shifting the rest of the data down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 442 (existing rows 442+ shift down to 444+)
$ws.Rows("442:443").Insert()

# ---- Row 442: Lane Late / Primera ----
$ws.Cells.Item(442, 1).Value = 3
$ws.Cells.Item(442, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(442, 3).Value = "Coquimbo"
$ws.Cells.Item(442, 4).Value = 44505
$ws.Cells.Item(442, 5).Value = 5
$ws.Cells.Item(442, 6).Value = "Fruta"
$ws.Cells.Item(442, 7).Value = 100102
$ws.Cells.Item(442, 8).Value = "Cítricos"
$ws.Cells.Item(442, 9).Value = 100102005
$ws.Cells.Item(442, 10).Value = "Naranja"
$ws.Cells.Item(442, 11).Value = "Lane Late"
$ws.Cells.Item(442, 12).Value = "Primera"
$ws.Cells.Item(442, 13).Value = 172
$ws.Cells.Item(442, 14).Value = 5000
$ws.Cells.Item(442, 15).Value = 6000
$ws.Cells.Item(442, 16).Value = 5494
$ws.Cells.Item(442, 17).Value = "`$/malla 13 kilos"
$ws.Cells.Item(442, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(442, 19).Value = 423
$ws.Cells.Item(442, 20).Value = 13

# ---- Row 443: Lane Late / Segunda ----
$ws.Cells.Item(443, 1).Value = 3
$ws.Cells.Item(443, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(443, 3).Value = "Coquimbo"
$ws.Cells.Item(443, 4).Value = 44505
$ws.Cells.Item(443, 5).Value = 5
$ws.Cells.Item(443, 6).Value = "Fruta"
$ws.Cells.Item(443, 7).Value = 100102
$ws.Cells.Item(443, 8).Value = "Cítricos"
$ws.Cells.Item(443, 9).Value = 100102005
$ws.Cells.Item(443, 10).Value = "Naranja"
$ws.Cells.Item(443, 11).Value = "Lane Late"
$ws.Cells.Item(443, 12).Value = "Segunda"
$ws.Cells.Item(443, 13).Value = 170
$ws.Cells.Item(443, 14).Value = 4000
$ws.Cells.Item(443, 15).Value = 4500
$ws.Cells.Item(443, 16).Value = 4265
$ws.Cells.Item(443, 17).Value = "`$/malla 13 kilos"
$ws.Cells.Item(443, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(443, 19).Value = 328
$ws.Cells.Item(443, 20).Value = 13

# ---- Append two more rows at the end (523 total rows) duplicating the
#      former last two data rows (old 520/521, now 522/523) ----

# Row 522: Lane Late / Segunda (copy of former row 520)
$ws.Cells.Item(522, 1).Value = 3
$ws.Cells.Item(522, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(522, 3).Value = "Coquimbo"
$ws.Cells.Item(522, 4).Value = 44491
$ws.Cells.Item(522, 5).Value = 5
$ws.Cells.Item(522, 6).Value = "Fruta"
$ws.Cells.Item(522, 7).Value = 100102
$ws.Cells.Item(522, 8).Value = "Cítricos"
$ws.Cells.Item(522, 9).Value = 100102005
$ws.Cells.Item(522, 10).Value = "Naranja"
$ws.Cells.Item(522, 11).Value = "Lane Late"
$ws.Cells.Item(522, 12).Value = "Segunda"
$ws.Cells.Item(522, 13).Value = 140
$ws.Cells.Item(522, 14).Value = 3500
$ws.Cells.Item(522, 15).Value = 4000
$ws.Cells.Item(522, 16).Value = 3750
$ws.Cells.Item(522, 17).Value = "`$/malla 13 kilos"
$ws.Cells.Item(522, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(522, 19).Value = 288
$ws.Cells.Item(522, 20).Value = 13

# Row 523: Navel Late / Primera (copy of former row 521)
$ws.Cells.Item(523, 1).Value = 3
$ws.Cells.Item(523, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(523, 3).Value = "Coquimbo"
$ws.Cells.Item(523, 4).Value = 44491
$ws.Cells.Item(523, 5).Value = 5
$ws.Cells.Item(523, 6).Value = "Fruta"
$ws.Cells.Item(523, 7).Value = 100102
$ws.Cells.Item(523, 8).Value = "Cítricos"
$ws.Cells.Item(523, 9).Value = 100102005
$ws.Cells.Item(523, 10).Value = "Naranja"
$ws.Cells.Item(523, 11).Value = "Navel Late"
$ws.Cells.Item(523, 12).Value = "Primera"
$ws.Cells.Item(523, 13).Value = 145
$ws.Cells.Item(523, 14).Value = 4500
$ws.Cells.Item(523, 15).Value = 5000
$ws.Cells.Item(523, 16).Value = 4759
$ws.Cells.Item(523, 17).Value = "`$/malla 13 kilos"
$ws.Cells.Item(523, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(523, 19).Value = 366
$ws.Cells.Item(523, 20).Value = 13

# Apply the same date number format used elsewhere in column D to the two
# newly inserted date cells (Insert() already copies formatting from the
# row above, but set it explicitly to be safe).
$ws.Cells.Item(442, 4).NumberFormat = $ws.Cells.Item(444, 4).NumberFormat
$ws.Cells.Item(443, 4).NumberFormat = $ws.Cells.Item(444, 4).NumberFormat
$ws.Cells.Item(522, 4).NumberFormat = $ws.Cells.Item(444, 4).NumberFormat
$ws.Cells.Item(523, 4).NumberFormat = $ws.Cells.Item(444, 4).NumberFormat
